$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -5
    "F4"  = -5
    "F6"  = -3
    "F7"  = -4
    "F8"  = -8
    "F11" = -5
    "F14" = -2
    "F16" = 3
    "F19" = -3
    "F21" = -7
    "F22" = -2
    "F23" = -3
    "F26" = 0
    "F32" = -2
    "F33" = -4
    "F39" = 5
    "F46" = -6
    "F48" = 2
    "F52" = -4
    "F60" = -4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
